# LV_Activities - 2nd May 2024
#
# Update the test-data workbook:
#  - Users sheet: replace tester name "James Craven" with "Indrajeet Singh"
#  - SaveActivityPopUpMsg sheet: replace error message text and make it the
#    active/selected sheet (moving focus away from ActivityStartDate)

$wb = $excel.ActiveWorkbook

# --- Users sheet -----------------------------------------------------
$wsUsers = $wb.Worksheets.Item("Users")
$wsUsers.Range("A2").Value = "Indrajeet Singh"
$wsUsers.Range("D6").Select()

# --- SaveActivityPopUpMsg sheet --------------------------------------
$wsSavePopup = $wb.Worksheets.Item("SaveActivityPopUpMsg")
$wsSavePopup.Range("A2").Value = "Complete this field."

# Make this sheet the active/selected tab (it becomes the workbook's
# active sheet, replacing ActivityStartDate).
$wsSavePopup.Activate()
$wsSavePopup.Range("B10").Select()
